$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current (before) list, rows 2-9:
# 2: Data drive lighting                                                    | 4
# 3: Make exporter a GUP.  Build UI & hook data into the max files          | 14
# 4: Replace XNAMATH with a Rorn maths library                              | 14
# 5: Refactoring - we need consistency across the board                     | 21
# 6: Error handling - go on a robustness run, set standards for future work | 21
# 7: Textured surfaces                                                      | 35
# 8: Vertex welding in model compiler                                      | 14
# 9: Complete the Rorn Maths library                                        | 35

# Target (after) list, rows 2-9:
# 2: Make exporter a GUP.  Build UI & hook data into the max files          | 14
# 3: Make the path from Max->Model Viewer seamless                          | 4
# 4: Replace XNAMATH with a Rorn maths library                              | 14
# 5: Refactoring - we need consistency across the board                     | 21
# 6: Error handling - go on a robustness run, set standards for future work | 21
# 7: Textured surfaces                                                      | 35
# 8: Vertex welding in model compiler                                      | 14
# 9: Complete the Rorn Maths library                                        | 35

$ws.Range("A2").Value = "Make exporter a GUP.  Build UI & hook data into the max files"
$ws.Range("B2").Value = 14

$ws.Range("A3").Value = "Make the path from Max->Model Viewer seamless"
$ws.Range("B3").Value = 4

$ws.Range("A4").Value = "Replace XNAMATH with a Rorn maths library"
$ws.Range("B4").Value = 14

$ws.Range("A5").Value = "Refactoring - we need consistency across the board"
$ws.Range("B5").Value = 21

$ws.Range("A6").Value = "Error handling - go on a robustness run, set standards for future work"
$ws.Range("B6").Value = 21

$ws.Range("A7").Value = "Textured surfaces"
$ws.Range("B7").Value = 35

$ws.Range("A8").Value = "Vertex welding in model compiler"
$ws.Range("B8").Value = 14

$ws.Range("A9").Value = "Complete the Rorn Maths library"
$ws.Range("B9").Value = 35

$ws.Range("B3").Select()
